# Add a new "Components Test" worksheet that builds a reusable subHeader
# component: the same masthead block used on "User Compliance Report"
# (logo/page, name/license block) plus three repeated "Annual Total: <range>"
# bands for 2015, 2016 and 2017 - one cycle-name + period combo per block.

$wb  = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("User Compliance Report")

# New sheet goes after the existing one.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Components Test"

# --- Rows 1-7: reuse the existing masthead layout (values + formats + merges) ---
$src.Range("A1:J7").Copy()
$ws.Range("A1:J7").PasteSpecial(-4104)   # xlPasteAll -> values + formats + merged cells
$src.Range("A1:J7").Copy()
$ws.Range("A1:J7").PasteSpecial(-4122)   # xlPasteFormats -> re-assert exact style ids

# --- Rows 8-11: two more "Annual Total" bands, same look as row 6/7 ---
$src.Range("A7:J7").Copy()
$ws.Range("A8:J8").PasteSpecial(-4122)
$ws.Range("A9:J9").PasteSpecial(-4122)
$ws.Range("A10:J10").PasteSpecial(-4122)
$ws.Range("A11:J11").PasteSpecial(-4122)

$ws.Range("A8:J9").Merge()
$ws.Range("A10:J11").Merge()

# Merging can nudge the covered cells onto a near-duplicate style; re-paste the
# clean row-7 format over the merged blocks so they land back on style id 4.
$src.Range("A7:J7").Copy()
$ws.Range("A8:J9").PasteSpecial(-4122)
$ws.Range("A10:J11").PasteSpecial(-4122)

# --- subHeader text: bold "Annual Total: " label + plain period ---
$ws.Range("A6").Value = "Annual Total: 1/1/2015 - 12/30/2015"
$ws.Range("A6").Characters(1, 14).Font.Bold = $true

$ws.Range("A8").Value = "Annual Total: 1/1/2016 - 12/30/2016"
$ws.Range("A8").Characters(1, 14).Font.Bold = $true

$ws.Range("A10").Value = "Annual Total: 1/1/2017 - 12/30/2017"
$ws.Range("A10").Characters(1, 14).Font.Bold = $true

# Keep the original report sheet as the active one.
$src.Activate()
